# Auto-generated: apply market-data refresh values per commit diff
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H43").Value = 2636
$ws.Range("J43").Value = 2657.923
$ws.Range("L43").Value = 2657.923
$ws.Range("N43").Value = -2795.923
$ws.Range("H98").Value = 985.52
$ws.Range("I98").Value = 770.95654
$ws.Range("J98").Value = 3453
$ws.Range("K98").Value = 770.95654
$ws.Range("L98").Value = 3453
$ws.Range("M98").Value = 727.04346
$ws.Range("N98").Value = -6449
$ws.Range("H122").Value = 985.52
$ws.Range("I122").Value = 770.95654
$ws.Range("J122").Value = 3453
$ws.Range("K122").Value = 2312.86962
$ws.Range("L122").Value = 10359
$ws.Range("M122").Value = 137.1303800000001
$ws.Range("N122").Value = -15259
$ws.Range("H129").Value = 962.6667
$ws.Range("I129").Value = 1033.8462
$ws.Range("J129").Value = 500
$ws.Range("K129").Value = 3101.5386
$ws.Range("L129").Value = 1500
$ws.Range("M129").Value = 1898.4614
$ws.Range("N129").Value = -11500
$ws.Range("H131").Value = 1197.5625
$ws.Range("I131").Value = 1197.5625
$ws.Range("K131").Value = 3592.6875
$ws.Range("M131").Value = 1447.3125
$ws.Range("H135").Value = 1971.9546
$ws.Range("I135").Value = 1232.25
$ws.Range("K135").Value = 11090.25
$ws.Range("M135").Value = -8555.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 816.8387
$ws.Range("I2").Value = 683.95654
$ws.Range("K2").Value = 683.95654
$ws.Range("M2").Value = -570.95654
$ws.Range("H32").Value = 31517.475
$ws.Range("I32").Value = 33284.234
$ws.Range("K32").Value = 33284.234
$ws.Range("M32").Value = -32997.234
$ws.Range("H116").Value = 816.8387
$ws.Range("I116").Value = 683.95654
$ws.Range("K116").Value = 683.95654
$ws.Range("M116").Value = 1610.04346
$ws.Range("H122").Value = 1999.3
$ws.Range("I122").Value = 1428.2858
$ws.Range("K122").Value = 4284.857400000001
$ws.Range("M122").Value = -1834.857400000001
$ws.Range("H132").Value = 103599.3
$ws.Range("I132").Value = 103599.3
$ws.Range("K132").Value = 310797.9
$ws.Range("M132").Value = -308267.9

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 816.8387
$ws.Range("I3").Value = 683.95654
$ws.Range("K3").Value = 683.95654
$ws.Range("M3").Value = -569.95654
$ws.Range("H75").Value = 25999
$ws.Range("J75").Value = 39999
$ws.Range("L75").Value = 39999
$ws.Range("N75").Value = -41871
$ws.Range("H78").Value = 25999
$ws.Range("J78").Value = 39999
$ws.Range("L78").Value = 119997
$ws.Range("N78").Value = -129357
$ws.Range("H86").Value = 2163.85
$ws.Range("I86").Value = 1950.4546
$ws.Range("J86").Value = 2424.6667
$ws.Range("K86").Value = 1950.4546
$ws.Range("L86").Value = 2424.6667
$ws.Range("M86").Value = -827.4546
$ws.Range("N86").Value = -4670.6667
$ws.Range("H89").Value = 2163.85
$ws.Range("I89").Value = 1950.4546
$ws.Range("J89").Value = 2424.6667
$ws.Range("K89").Value = 9752.273000000001
$ws.Range("L89").Value = 12123.3335
$ws.Range("M89").Value = -4136.273000000001
$ws.Range("N89").Value = -23355.3335
$ws.Range("H134").Value = 2516.7
$ws.Range("I134").Value = 2516.7
$ws.Range("K134").Value = 7550.099999999999
$ws.Range("M134").Value = -5015.099999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3096.7097
$ws.Range("I31").Value = 2426.9167
$ws.Range("J31").Value = 5393.143
$ws.Range("K31").Value = 2426.9167
$ws.Range("L31").Value = 5393.143
$ws.Range("M31").Value = -2131.9167
$ws.Range("N31").Value = -5983.143
$ws.Range("H34").Value = 3096.7097
$ws.Range("I34").Value = 2426.9167
$ws.Range("J34").Value = 5393.143
$ws.Range("K34").Value = 2426.9167
$ws.Range("L34").Value = 5393.143
$ws.Range("M34").Value = -2224.9167
$ws.Range("N34").Value = -5797.143
$ws.Range("H86").Value = 35281.92
$ws.Range("I86").Value = 48741.332
$ws.Range("J86").Value = 4998.25
$ws.Range("K86").Value = 48741.332
$ws.Range("L86").Value = 4998.25
$ws.Range("M86").Value = -47618.332
$ws.Range("N86").Value = -7244.25
$ws.Range("H89").Value = 35281.92
$ws.Range("I89").Value = 48741.332
$ws.Range("J89").Value = 4998.25
$ws.Range("K89").Value = 243706.66
$ws.Range("L89").Value = 24991.25
$ws.Range("M89").Value = -238090.66
$ws.Range("N89").Value = -36223.25
$ws.Range("H122").Value = 1793.3077
$ws.Range("I122").Value = 1813.6666
$ws.Range("K122").Value = 5440.9998
$ws.Range("M122").Value = -2990.9998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 333.33334
$ws.Range("I11").Value = 333.33334
$ws.Range("K11").Value = 1000.00002
$ws.Range("M11").Value = -860.0000200000001
$ws.Range("H26").Value = 230.2
$ws.Range("I26").Value = 200.33333
$ws.Range("K26").Value = 600.99999
$ws.Range("M26").Value = -312.99999
$ws.Range("H93").Value = 8124.6
$ws.Range("J93").Value = 10000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -33744
$ws.Range("H109").Value = 3191.5
$ws.Range("I109").Value = 2787.25
$ws.Range("K109").Value = 8361.75
$ws.Range("M109").Value = -7321.75
$ws.Range("H131").Value = 2133253.8
$ws.Range("I131").Value = 2559.6667
$ws.Range("J131").Value = 2445062.5
$ws.Range("K131").Value = 7679.000100000001
$ws.Range("L131").Value = 7335187.5
$ws.Range("M131").Value = -2639.000100000001
$ws.Range("N131").Value = -7345267.5
$ws.Range("H139").Value = 2772.875
$ws.Range("I139").Value = 2534.1428
$ws.Range("K139").Value = 7602.428400000001
$ws.Range("M139").Value = -2462.428400000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2002
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -10008
$ws.Range("N83").ClearContents()
$ws.Range("H107").Value = 42566.832
$ws.Range("I107").Value = 56436.723
$ws.Range("J107").Value = 957.1667
$ws.Range("K107").Value = 56436.723
$ws.Range("L107").Value = 957.1667
$ws.Range("M107").Value = -54516.723
$ws.Range("N107").Value = -4797.1667
$ws.Range("H132").Value = 54668.05
$ws.Range("I132").Value = 78976.38
$ws.Range("K132").Value = 236929.14
$ws.Range("M132").Value = -234399.14

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 598.4
$ws.Range("I16").Value = 657.55554
$ws.Range("J16").Value = 66
$ws.Range("K16").Value = 657.55554
$ws.Range("L16").Value = 66
$ws.Range("M16").Value = -487.55554
$ws.Range("N16").Value = -406
$ws.Range("H22").Value = 2399.1765
$ws.Range("I22").Value = 619.2
$ws.Range("J22").Value = 4942
$ws.Range("K22").Value = 619.2
$ws.Range("L22").Value = 4942
$ws.Range("M22").Value = -324.2
$ws.Range("N22").Value = -5532
$ws.Range("H27").Value = 2399.1765
$ws.Range("I27").Value = 619.2
$ws.Range("J27").Value = 4942
$ws.Range("K27").Value = 619.2
$ws.Range("L27").Value = 4942
$ws.Range("M27").Value = -512.2
$ws.Range("N27").Value = -5156
$ws.Range("H61").Value = 4579.8
$ws.Range("I61").Value = 4579.8
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4579.8
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4377.8
$ws.Range("N61").ClearContents()
$ws.Range("H82").Value = 3009.6667
$ws.Range("J82").Value = 3085.05
$ws.Range("L82").Value = 3085.05
$ws.Range("N82").Value = -3807.05
$ws.Range("H85").Value = 3009.6667
$ws.Range("J85").Value = 3085.05
$ws.Range("L85").Value = 3085.05
$ws.Range("N85").Value = -5581.05
$ws.Range("H113").Value = 4579.8
$ws.Range("I113").Value = 4579.8
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4579.8
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2409.8
$ws.Range("N113").ClearContents()
$ws.Range("H136").Value = 4757.5
$ws.Range("I136").Value = 3172.5
$ws.Range("K136").Value = 9517.5
$ws.Range("M136").Value = -6967.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 76479.34
$ws.Range("I126").Value = 92833.81
$ws.Range("J126").Value = 7790.6
$ws.Range("K126").Value = 278501.43
$ws.Range("L126").Value = 23371.8
$ws.Range("M126").Value = -276031.43
$ws.Range("N126").Value = -28311.8
$ws.Range("H136").Value = 2096.2563
$ws.Range("I136").Value = 1763.3125
$ws.Range("J136").Value = 3618.2856
$ws.Range("K136").Value = 5289.9375
$ws.Range("L136").Value = 10854.8568
$ws.Range("M136").Value = -2739.9375
